# Chapter 9 lecture edit
# 1) Refresh the cached "datetimeFigureOut" date field text (slide master +
#    every slide layout) from 7/16/18 -> 11/27/18.
# 2) Re-title slide 1: "Chapter 3" -> "Chapter 9" (last character edited in
#    place, which is how PowerPoint splits the run into an untouched
#    "Chapter " run and a freshly-edited "9" run).
# 3) Collapse the two-line subtitle into a single line: "Spatial Statistics".

$p = $ppt.ActivePresentation

$oldDate = "7/16/18"
$newDate = "11/27/18"

# --- Slide master date placeholder ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every slide layout's date placeholder ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide 1 title: "Chapter 3" -> "Chapter 9" ---
$slide = $p.Slides.Item(1)
$title = $slide.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$lastChar = $titleRange.Characters($titleRange.Length, 1)
$lastChar.Text = "9"

# --- Slide 1 subtitle: two lines -> single "Spatial Statistics" line ---
$subtitle = $slide.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = "Spatial Statistics"
